# Update market-price / profit figures across several Leve-profit sheets.
# Values below come from a refreshed Universalis market-data pull; only the
# price/profit columns (H:N) change, row-by-row, per worksheet.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 7
$ws.Range("H7").Value = 250
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 250
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 250
$ws.Range("M7").ClearContents()
$ws.Range("N7").Value = -474
# Row 11
$ws.Range("H11").Value = 99
$ws.Range("I11").Value = 99
$ws.Range("K11").Value = 99
$ws.Range("M11").Value = 41
# Row 14
$ws.Range("H14").Value = 250
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 250
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 250
$ws.Range("M14").ClearContents()
$ws.Range("N14").Value = -632
# Row 18
$ws.Range("H18").Value = 4249
$ws.Range("J18").Value = 3500
$ws.Range("L18").Value = 3500
$ws.Range("N18").Value = -4068
# Row 21
$ws.Range("H21").Value = 808.3333
$ws.Range("J21").Value = 750
$ws.Range("L21").Value = 750
$ws.Range("N21").Value = -1686
# Row 23
$ws.Range("H23").Value = 808.3333
$ws.Range("J23").Value = 750
$ws.Range("L23").Value = 750
$ws.Range("N23").Value = -1218
# Row 38
$ws.Range("H38").Value = 56.666668
$ws.Range("I38").Value = 56.666668
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 170.000004
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = 201.999996
$ws.Range("N38").ClearContents()
# Row 55
$ws.Range("H55").Value = 499.25
$ws.Range("I55").Value = 548.5
$ws.Range("K55").Value = 548.5
$ws.Range("M55").Value = -334.5
# Row 107
$ws.Range("H107").Value = 0
$ws.Range("I107").Value = 0
$ws.Range("K107").Value = 0
$ws.Range("M107").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
# Row 5
$ws.Range("H5").Value = 112.5
$ws.Range("I5").Value = 100
$ws.Range("J5").Value = 125
$ws.Range("K5").Value = 100
$ws.Range("L5").Value = 125
$ws.Range("M5").Value = 12
$ws.Range("N5").Value = -349
# Row 22
$ws.Range("H22").Value = 2000
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()
# Row 24
$ws.Range("H24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("N24").ClearContents()
# Row 88
$ws.Range("H88").Value = 2580.8333
$ws.Range("J88").Value = 2000
$ws.Range("L88").Value = 2000
$ws.Range("N88").Value = -2812
# Row 91
$ws.Range("H91").Value = 2580.8333
$ws.Range("J91").Value = 2000
$ws.Range("L91").Value = 2000
$ws.Range("N91").Value = -4808
# Row 92
$ws.Range("H92").Value = 41274
$ws.Range("J92").Value = 41274
$ws.Range("L92").Value = 41274
$ws.Range("N92").Value = -46266
# Row 97
$ws.Range("H97").Value = 607.25
$ws.Range("I97").Value = 607.25
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 607.25
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -111.25
$ws.Range("N97").ClearContents()
# Row 100
$ws.Range("H100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
# Row 35
$ws.Range("H35").Value = 4459.5
$ws.Range("I35").Value = 950
$ws.Range("K35").Value = 950
$ws.Range("M35").Value = -656
# Row 38
$ws.Range("H38").Value = 9679.333000000001
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents()
# Row 46
$ws.Range("H46").Value = 9679.333000000001
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()
# Row 62
$ws.Range("H62").Value = 566.6667
$ws.Range("J62").Value = 700
$ws.Range("L62").Value = 700
$ws.Range("N62").Value = -1948
# Row 65
$ws.Range("H65").Value = 566.6667
$ws.Range("J65").Value = 700
$ws.Range("L65").Value = 3500
$ws.Range("N65").Value = -9740
# Row 105
$ws.Range("H105").Value = 2828.2856
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()
# Row 134
$ws.Range("H134").Value = 6529.6665
$ws.Range("J134").Value = 17000
$ws.Range("L134").Value = 51000
$ws.Range("N134").Value = -56070

$ws = $wb.Worksheets.Item("CUL")
# Row 13
$ws.Range("H13").Value = 4133.3335
$ws.Range("I13").Value = 4200
$ws.Range("J13").Value = 4000
$ws.Range("K13").Value = 12600
$ws.Range("L13").Value = 12000
$ws.Range("M13").Value = -12432
$ws.Range("N13").Value = -12336
# Row 17
$ws.Range("H17").Value = 257.25
$ws.Range("J17").Value = 334.66666
$ws.Range("L17").Value = 1003.99998
$ws.Range("N17").Value = -1341.99998
# Row 42
$ws.Range("H42").Value = 1000
$ws.Range("J42").Value = 1000
$ws.Range("L42").Value = 3000
$ws.Range("N42").Value = -4068
# Row 63
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("M63").ClearContents()
# Row 66
$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("M66").ClearContents()
# Row 141
$ws.Range("H141").Value = 2894.5
$ws.Range("I141").Value = 2894.5
$ws.Range("K141").Value = 8683.5
$ws.Range("M141").Value = -3503.5

$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 2576.75
$ws.Range("J80").Value = 2576.75
$ws.Range("L80").Value = 2576.75
$ws.Range("N80").Value = -4572.75
# Row 83
$ws.Range("H83").Value = 2576.75
$ws.Range("J83").Value = 2576.75
$ws.Range("L83").Value = 12883.75
$ws.Range("N83").Value = -22867.75

$ws = $wb.Worksheets.Item("LTW")
# Row 4
$ws.Range("H4").Value = 16600
# Row 22
$ws.Range("H22").Value = 2308
$ws.Range("J22").Value = 1893.25
$ws.Range("L22").Value = 1893.25
$ws.Range("N22").Value = -2483.25
# Row 27
$ws.Range("H27").Value = 2308
$ws.Range("J27").Value = 1893.25
$ws.Range("L27").Value = 1893.25
$ws.Range("N27").Value = -2107.25
# Row 28
$ws.Range("H28").Value = 16600
# Row 37
$ws.Range("H37").Value = 16600
# Row 43
$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()
# Row 46
$ws.Range("H46").Value = 5262.5
$ws.Range("J46").Value = 3250
$ws.Range("L46").Value = 3250
$ws.Range("N46").Value = -3626

$ws = $wb.Worksheets.Item("WVR")
# Row 18
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()
# Row 62
$ws.Range("H62").Value = 54325
$ws.Range("J62").Value = 54325
$ws.Range("L62").Value = 54325
$ws.Range("N62").Value = -55573
# Row 63
$ws.Range("H63").Value = 26309.75
$ws.Range("J63").Value = 31748
$ws.Range("L63").Value = 31748
$ws.Range("N63").Value = -32996
# Row 65
$ws.Range("H65").Value = 54325
$ws.Range("J65").Value = 54325
$ws.Range("L65").Value = 271625
$ws.Range("N65").Value = -277865
# Row 66
$ws.Range("H66").Value = 26309.75
$ws.Range("J66").Value = 31748
$ws.Range("L66").Value = 95244
$ws.Range("N66").Value = -101484
